# Append a new paragraph at the very end of the document (after the
# "...New table" paragraph, before the sectPr) containing the Chap 7 notes.
$d = $word.ActiveDocument

# Collapse a range at the end of the document content and insert a new
# paragraph mark there; the new (empty) paragraph inherits the run/paragraph
# formatting (Times New Roman) from the preceding paragraph, same as Word
# does when you place the cursor at the end of the doc and hit Enter.
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

# The newly created paragraph is now the last paragraph in the document;
# fill in its text.
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "Chap 7: A master refactor changing from Identity User to Employee just to satisfy a single mapping line code on Leave Allocation control. Also this chapter allows us to learn about refactoring in C# when business requirement"
